$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.85
$ws.Range("I2").Value = 4.3
$ws.Range("J2").Value = 2.45
$ws.Range("L2").Value = 4.6
$ws.Range("N2").Value = 6.3
$ws.Range("P2").Value = 2.72
$ws.Range("Q2").Value = 2.18
$ws.Range("T2").Value = 2.57
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.72
$ws.Range("W2").Value = 5.8
$ws.Range("X2").Value = 7.8
$ws.Range("Y2").Value = 8.5
$ws.Range("AA2").Value = 17
$ws.Range("AB2").Value = 35
$ws.Range("AC2").Value = 6.3
$ws.Range("AH2").Value = 10.25
$ws.Range("AJ2").Value = 14
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 45
$ws.Range("AM2").Value = 55
$ws.Range("AO2").Value = 9.5
$ws.Range("AP2").Value = 20
$ws.Range("AQ2").Value = 35
$ws.Range("AR2").Value = 75
$ws.Range("AT2").Value = 2.57
$ws.Range("AW2").Value = 5.9
$ws.Range("AX2").Value = 25
$ws.Range("AY2").Value = 32
$ws.Range("AZ2").Value = 150
$ws.Range("BA2").Value = 175
$ws.Range("BB2").Value = 450

# Row 4
$ws.Range("G4").Value = 1.53
$ws.Range("I4").Value = 5.5
$ws.Range("L4").Value = 5
$ws.Range("U4").Value = 1.57
$ws.Range("V4").Value = 2.25
$ws.Range("AM4").Value = 34
$ws.Range("AY4").Value = 26

# Row 5
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 19
$ws.Range("Q5").Value = 1.5
$ws.Range("R5").Value = 2.5

# Row 6
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.4
